$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Clcf1"
$ws.Cells.Item(2,3).Value = "Lifr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.9214586666666666
$ws.Cells.Item(2,8).Value = 2.764376
$ws.Cells.Item(2,9).Value = 0.08041853843186561
$ws.Cells.Item(2,10).Value = 0.08041853843186561
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 35.337883
$ws.Cells.Item(2,14).Value = 106.013649
$ws.Cells.Item(2,15).Value = 0.3968231145247413
$ws.Cells.Item(2,16).Value = 0.3968231145247413
$ws.Cells.Item(2,17).Value = 32.56239855200266
$ws.Cells.Item(2,18).Value = 293.061586968024
$ws.Cells.Item(2,19).Value = 0.03191193488606052
$ws.Cells.Item(2,20).Value = 0.03191193488606051

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Clcf1"
$ws.Cells.Item(3,3).Value = "Lifr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.9214586666666666
$ws.Cells.Item(3,8).Value = 2.764376
$ws.Cells.Item(3,9).Value = 0.08041853843186561
$ws.Cells.Item(3,10).Value = 0.08041853843186561
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 40.49537033333333
$ws.Cells.Item(3,14).Value = 121.486111
$ws.Cells.Item(3,15).Value = 0.4547385869013757
$ws.Cells.Item(3,16).Value = 0.4547385869013756
$ws.Cells.Item(3,17).Value = 37.31480995352622
$ws.Cells.Item(3,18).Value = 335.833289581736
$ws.Cells.Item(3,19).Value = 0.03656941252718054
$ws.Cells.Item(3,20).Value = 0.03656941252718053

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Clcf1"
$ws.Cells.Item(4,3).Value = "Lifr"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.9214586666666666
$ws.Cells.Item(4,8).Value = 2.764376
$ws.Cells.Item(4,9).Value = 0.08041853843186561
$ws.Cells.Item(4,10).Value = 0.08041853843186561
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 13.218724
$ws.Cells.Item(4,14).Value = 39.656172
$ws.Cells.Item(4,15).Value = 0.1484382985738831
$ws.Cells.Item(4,16).Value = 0.148438298573883
$ws.Cells.Item(4,17).Value = 12.18050779207467
$ws.Cells.Item(4,18).Value = 109.624570128672
$ws.Cells.Item(4,19).Value = 0.01193719101862456
$ws.Cells.Item(4,20).Value = 0.01193719101862456

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Clcf1"
$ws.Cells.Item(5,3).Value = "Lifr"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.913147
$ws.Cells.Item(5,8).Value = 5.739441
$ws.Cells.Item(5,9).Value = 0.1669662363715809
$ws.Cells.Item(5,10).Value = 0.1669662363715809
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 35.337883
$ws.Cells.Item(5,14).Value = 106.013649
$ws.Cells.Item(5,15).Value = 0.3968231145247413
$ws.Cells.Item(5,16).Value = 0.3968231145247413
$ws.Cells.Item(5,17).Value = 67.606564847801
$ws.Cells.Item(5,18).Value = 608.4590836302091
$ws.Cells.Item(5,19).Value = 0.06625606193744486
$ws.Cells.Item(5,20).Value = 0.06625606193744486

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Clcf1"
$ws.Cells.Item(6,3).Value = "Lifr"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.913147
$ws.Cells.Item(6,8).Value = 5.739441
$ws.Cells.Item(6,9).Value = 0.1669662363715809
$ws.Cells.Item(6,10).Value = 0.1669662363715809
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 40.49537033333333
$ws.Cells.Item(6,14).Value = 121.486111
$ws.Cells.Item(6,15).Value = 0.4547385869013757
$ws.Cells.Item(6,16).Value = 0.4547385869013756
$ws.Cells.Item(6,17).Value = 77.47359626710568
$ws.Cells.Item(6,18).Value = 697.262366403951
$ws.Cells.Item(6,19).Value = 0.07592599038785376
$ws.Cells.Item(6,20).Value = 0.07592599038785375

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Clcf1"
$ws.Cells.Item(7,3).Value = "Lifr"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.913147
$ws.Cells.Item(7,8).Value = 5.739441
$ws.Cells.Item(7,9).Value = 0.1669662363715809
$ws.Cells.Item(7,10).Value = 0.1669662363715809
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 13.218724
$ws.Cells.Item(7,14).Value = 39.656172
$ws.Cells.Item(7,15).Value = 0.1484382985738831
$ws.Cells.Item(7,16).Value = 0.148438298573883
$ws.Cells.Item(7,17).Value = 25.289362164428
$ws.Cells.Item(7,18).Value = 227.604259479852
$ws.Cells.Item(7,19).Value = 0.02478418404628226
$ws.Cells.Item(7,20).Value = 0.02478418404628225

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Clcf1"
$ws.Cells.Item(8,3).Value = "Lifr"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 8.623680999999999
$ws.Cells.Item(8,8).Value = 25.871043
$ws.Cells.Item(8,9).Value = 0.7526152251965536
$ws.Cells.Item(8,10).Value = 0.7526152251965536
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 35.337883
$ws.Cells.Item(8,14).Value = 106.013649
$ws.Cells.Item(8,15).Value = 0.3968231145247413
$ws.Cells.Item(8,16).Value = 0.3968231145247413
$ws.Cells.Item(8,17).Value = 304.742630207323
$ws.Cells.Item(8,18).Value = 2742.683671865907
$ws.Cells.Item(8,19).Value = 0.298655117701236
$ws.Cells.Item(8,20).Value = 0.2986551177012359

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Clcf1"
$ws.Cells.Item(9,3).Value = "Lifr"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 8.623680999999999
$ws.Cells.Item(9,8).Value = 25.871043
$ws.Cells.Item(9,9).Value = 0.7526152251965536
$ws.Cells.Item(9,10).Value = 0.7526152251965536
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 40.49537033333333
$ws.Cells.Item(9,14).Value = 121.486111
$ws.Cells.Item(9,15).Value = 0.4547385869013757
$ws.Cells.Item(9,16).Value = 0.4547385869013756
$ws.Cells.Item(9,17).Value = 349.2191557315303
$ws.Cells.Item(9,18).Value = 3142.972401583773
$ws.Cells.Item(9,19).Value = 0.3422431839863414
$ws.Cells.Item(9,20).Value = 0.3422431839863413

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Clcf1"
$ws.Cells.Item(10,3).Value = "Lifr"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 8.623680999999999
$ws.Cells.Item(10,8).Value = 25.871043
$ws.Cells.Item(10,9).Value = 0.7526152251965536
$ws.Cells.Item(10,10).Value = 0.7526152251965536
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 13.218724
$ws.Cells.Item(10,14).Value = 39.656172
$ws.Cells.Item(10,15).Value = 0.1484382985738831
$ws.Cells.Item(10,16).Value = 0.148438298573883
$ws.Cells.Item(10,17).Value = 113.994059003044
$ws.Cells.Item(10,18).Value = 1025.946531027396
$ws.Cells.Item(10,19).Value = 0.1117169235089763
$ws.Cells.Item(10,20).Value = 0.1117169235089762
